# Horarios actualizados Línea 141 - 723
# Updates the "Última actualización" timestamp, "Total filas" counts, a couple
# of existing rows whose elapsed-minutes got recalculated against the new
# scrape time, and appends newly scraped arrival rows to the LP1912 and
# LP1912-215 sheets.

$wb = $excel.ActiveWorkbook

$oldTime = "02:49:45"
$newTime = "03:00:18"

# ---------------------------------------------------------------------
# Sheet "LP1912"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: $newTime"
$ws1.Range("A3").Value = "Total filas: 12"

# Existing row 11 and row 14 get refreshed Hora_Scrap + Minutos values
$ws1.Cells.Item(11, 1).Value = $newTime
$ws1.Cells.Item(11, 4).Value = 48

$ws1.Cells.Item(14, 1).Value = $newTime
$ws1.Cells.Item(14, 4).Value = 61

# New rows 16 and 17
$ws1.Cells.Item(16, 1).Value = $newTime
$ws1.Cells.Item(16, 2).Value = "04:44"
$ws1.Cells.Item(16, 3).Value = "215_ALUAR"
$ws1.Cells.Item(16, 4).Value = 104
$ws1.Cells.Item(16, 5).Value = "LP1912"

$ws1.Cells.Item(17, 1).Value = $newTime
$ws1.Cells.Item(17, 2).Value = "04:53"
$ws1.Cells.Item(17, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(17, 4).Value = 113
$ws1.Cells.Item(17, 5).Value = "LP1912"

# ---------------------------------------------------------------------
# Sheet "LP1912-215"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: $newTime"
$ws2.Range("A3").Value = "Total filas: 7"

# New row 12
$ws2.Cells.Item(12, 1).Value = $newTime
$ws2.Cells.Item(12, 2).Value = "04:44"
$ws2.Cells.Item(12, 3).Value = "215_ALUAR"
$ws2.Cells.Item(12, 4).Value = 104
$ws2.Cells.Item(12, 5).Value = "LP1912"

# ---------------------------------------------------------------------
# Sheet "6203-6173"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: $newTime"
